# 20-Feb-2017-2045 - Committing Changes of Client User along with the xlsx and .Java Files
#
# Populates the "Results" column for the existing Test Cases / CRCreation
# rows that were run:
#   - Test Cases!D3  -> PASS
#   - Test Cases!D4  -> PASS
#   - CRCreation!AU2 -> SKIP
#   - CRCreation!AU3 -> PASS

$wb = $excel.ActiveWorkbook

# NOTE: "SKIP" must be written before "PASS" so the two new shared-string
# table entries land in the same order as the target workbook (SKIP, PASS).
$wsCRCreation = $wb.Worksheets.Item("CRCreation")
$wsCRCreation.Range("AU2").Value = "SKIP"

$wsTestCases = $wb.Worksheets.Item("Test Cases")
$wsTestCases.Range("D3").Value = "PASS"
$wsTestCases.Range("D4").Value = "PASS"

$wsCRCreation.Range("AU3").Value = "PASS"
